$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are forced to text (leading-apostrophe trick, like typing
# into the Excel UI) because several values are plain numerals (e.g. "1.00", "605.65")
# that Excel would otherwise silently convert to numbers, dropping the original
# formatting. Style is reset to "Normal" right after so no stray number-format style
# sticks to the cell.

$ws.Range("D2").Value = "`'72.361.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.47%  "
$ws.Range("D3").Value = "`'2.632.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "`'605.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "`'179.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "`'0.528"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("E9").Value = "  +8.93%  "
$ws.Range("D10").Value = "`'2.630.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("D12").Value = "`'0.355"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "`'0.0000190"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.41%  "
$ws.Range("D15").Value = "`'3.106.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "`'26.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "`'72.281.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("D18").Value = "`'2.636.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("D19").Value = "`'11.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.46%  "
$ws.Range("D20").Value = "`'382.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.62%  "
$ws.Range("D21").Value = "`'7.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  +18.08%  "
$ws.Range("D24").Value = "`'73.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.32%  "
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").Value = "`'4.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "`'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "`'10.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.31%  "
$ws.Range("D28").Value = "`'2.766.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "`'0.0₃0961"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.69%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "`'8.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.69%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "`'519.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "`'165.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("D37").Value = "`'19.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("D39").Value = "`'19.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").Value = "`'0.112"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.60%  "
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("D42").Value = "`'5.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "`'2.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("D46").Value = "`'39.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "`'150.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("D49").Value = "`'0.547"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("E50").Value = "  +4.18%  "
$ws.Range("D51").Value = "`'0.0₆0265"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.73%  "
